$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 97, shifting existing rows 97-232 down to 98-233.
$ws.Rows("97").Insert()

# Populate the newly inserted row 97 with this week's record for
# Macroferia Regional de Talca - Ajo - Chino - Primera.
$ws.Range("A97").Value = 5
$ws.Range("B97").Value = "Macroferia Regional de Talca"
$ws.Range("C97").Value = "Maule"
$ws.Range("D97").Value = 44579
$ws.Range("E97").Value = 7
$ws.Range("F97").Value = 100112003
$ws.Range("G97").Value = "Ajo"
$ws.Range("H97").Value = "Chino"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 150
$ws.Range("K97").Value = 19000
$ws.Range("L97").Value = 20000
$ws.Range("M97").Value = 19667
$ws.Range("N97").Value = "$/caja 10 kilos"
$ws.Range("O97").Value = "China"
$ws.Range("P97").Value = 1967
$ws.Range("Q97").Value = 10
$ws.Range("R97").Value = "Hortaliza"
